# "New slides, new scripts" — append a new "Vocabulary" slide at the end
# of the deck (position 5), using the same title+body(idx=10) layout that
# the other content slides (slide2/3/4) already use.

$p = $ppt.ActivePresentation

# ppLayout=3 -> the custom "Classic slide, one column" layout (the 3rd
# slide layout in this deck), which exposes a `title` placeholder and a
# `body` placeholder at idx=10 — exactly the shapes the new slide needs.
$s = $p.Slides.Add($p.Slides.Count + 1, 3)

# --- Title placeholder -----------------------------------------------
$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.Text = "Vocabulary"

# --- Body placeholder --------------------------------------------------
$body = $s.Shapes.Item(2).TextFrame.TextRange

# Seed with a single-paragraph placeholder first so PowerPoint stamps the
# `lang` attribute onto every run once the real, multi-paragraph text is
# written in afterwards.
$body.Text = "placeholder"

$para1 = "Vocabulary: limited set of discrete items"
$para2 = "For example: all words* in the English language"
$para3 = "Vocabulary is used for both inputs and outputs"
$para4 = "x"
$para5 = "What is the vocabulary of a Transformer generating text?"

$body.Text = $para1 + "`r" + $para2 + "`r" + $para3 + "`r" + $para4 + "`r" + $para5

# Second-level bullets for paragraphs 2 and 3 (0-based lvl="1").
$body.Paragraphs(2, 1).IndentLevel = 2
$body.Paragraphs(3, 1).IndentLevel = 2

# Clear the 4th paragraph (seeded with a placeholder "x" above) back down
# to a genuinely empty line/paragraph mark.
$body.Paragraphs(4, 1).Text = ""

# Bold spans inside paragraph 1: "limited set" and "discrete items".
$body.Characters(13, 11).Font.Bold = $true
$body.Characters(28, 14).Font.Bold = $true

# Bold spans inside paragraph 3: "inputs" and "outputs".
$p3start = ($para1.Length + 1) + ($para2.Length + 1)
$body.Characters($p3start + 29, 6).Font.Bold = $true
$body.Characters($p3start + 40, 7).Font.Bold = $true
